# Update the "Pais" COVID snapshot: refresh the timestamp, update the
# per-country case/recovered/death counters for the countries whose
# figures moved, and fix the country-name ordering for the rows whose
# rank changed as a result (Bielorrusia/Catar, Australia/Kuwait,
# Crucero/Uruguay/Mali/Maldivas, Seychelles/Montserrat).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 15:04"
$ws.Range("B4").Value = 1292917
$ws.Range("C4").Value = 294
$ws.Range("E4").Value = 998724
$ws.Range("B19").Value = 42093
$ws.Range("C19").Value = 319
$ws.Range("E19").Value = 36484
$ws.Range("G19").Value = 71
$ws.Range("H19").Value = 5359
$ws.Range("B20").Value = 35432
$ws.Range("C20").Value = 1701
$ws.Range("D20").Value = 9120
$ws.Range("E20").Value = 26083
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 229
$ws.Range("B24").Value = 27268
$ws.Range("C24").Value = 553
$ws.Range("D24").Value = 2422
$ws.Range("E24").Value = 23732
$ws.Range("F24").Value = 127
$ws.Range("G24").Value = 9
$ws.Range("H24").Value = 1114
$ws.Range("B25").Value = 26435
$ws.Range("C25").Value = 1791
$ws.Range("E25").Value = 18306
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 599
$ws.Range("B26").Value = 25265
$ws.Range("C26").Value = 642
$ws.Range("E26").Value = 17119
$ws.Range("G26").Value = 135
$ws.Range("H26").Value = 3175
$ws.Range("A30").Value = "Catar"
$ws.Range("B30").Value = 20201
$ws.Range("C30").Value = 1311
$ws.Range("D30").Value = 2370
$ws.Range("E30").Value = 17819
$ws.Range("F30").Value = 72
$ws.Range("H30").Value = 12
$ws.Range("A31").Value = "Bielorrusia"
$ws.Range("B31").Value = 20168
$ws.Range("D31").Value = 5067
$ws.Range("E31").Value = 14985
$ws.Range("F31").Value = 92
$ws.Range("H31").Value = 116
$ws.Range("B43").Value = 10218
$ws.Range("C43").Value = 135
$ws.Range("D43").Value = 7927
$ws.Range("E43").Value = 1769
$ws.Range("G43").Value = 8
$ws.Range("H43").Value = 522
$ws.Range("A52").Value = "Kuwait"
$ws.Range("B52").Value = 7208
$ws.Range("C52").Value = 641
$ws.Range("D52").Value = 2466
$ws.Range("E52").Value = 4695
$ws.Range("F52").Value = 91
$ws.Range("G52").Value = 3
$ws.Range("H52").Value = 47
$ws.Range("A53").Value = "Australia"
$ws.Range("B53").Value = 6914
$ws.Range("C53").Value = 18
$ws.Range("D53").Value = 6079
$ws.Range("E53").Value = 738
$ws.Range("F53").Value = 21
$ws.Range("H53").Value = 97
$ws.Range("B55").Value = 5738
$ws.Range("C55").Value = 65
$ws.Range("E55").Value = 1978
$ws.Range("F55").Value = 45
$ws.Range("G55").Value = 5
$ws.Range("H55").Value = 260
$ws.Range("B59").Value = 4794
$ws.Range("C59").Value = 216
$ws.Range("E59").Value = 3181
$ws.Range("B61").Value = 4404
$ws.Range("C61").Value = 205
$ws.Range("D61").Value = 2027
$ws.Range("E61").Value = 2369
$ws.Range("F61").Value = 2
$ws.Range("B72").Value = 2324
$ws.Range("C72").Value = 26
$ws.Range("D72").Value = 1748
$ws.Range("E72").Value = 566
$ws.Range("B75").Value = 2161
$ws.Range("C75").Value = 36
$ws.Range("D75").Value = 1689
$ws.Range("E75").Value = 386
$ws.Range("F75").Value = 13
$ws.Range("A110").Value = "Maldivas"
$ws.Range("B110").Value = 734
$ws.Range("C110").Value = 86
$ws.Range("D110").Value = 20
$ws.Range("E110").Value = 711
$ws.Range("F110").Value = 2
$ws.Range("H110").Value = 3
$ws.Range("A111").Value = "Crucero"
$ws.Range("B111").Value = 712
$ws.Range("D111").Value = 645
$ws.Range("E111").Value = 54
$ws.Range("F111").Value = 4
$ws.Range("H111").Value = 13
$ws.Range("A112").Value = "Uruguay"
$ws.Range("B112").Value = 684
$ws.Range("D112").Value = 492
$ws.Range("E112").Value = 175
$ws.Range("F112").Value = 9
$ws.Range("H112").Value = 17
$ws.Range("A113").Value = "Mali"
$ws.Range("B113").Value = 650
$ws.Range("D113").Value = 271
$ws.Range("E113").Value = 347
$ws.Range("F113").Value = 0
$ws.Range("H113").Value = 32
$ws.Range("A205").Value = "Montserrat"
$ws.Range("D205").Value = 7
$ws.Range("F205").Value = 1
$ws.Range("H205").Value = 1
$ws.Range("A206").Value = "Seychelles"
$ws.Range("D206").Value = 8
$ws.Range("F206").Value = 0
$ws.Range("H206").Value = 0
